$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("for each of the last ", $true, $false, $false, $false, $false,
              $true, 1, $false, "in each of the last ", 2)
